# Update "想去人数" (interested-count) values on both the "展览" sheet
# and the "全部类型" sheet (which mirrors the same rows), matching the
# site's regenerated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 580
    $ws.Range("F3").Value = 3638
    $ws.Range("F5").Value = 703
}
